# Add a "constraint" column to the survey sheet and a new "age" (decimal)
# question row that uses it, per the XLSForm-style layout already in the
# workbook:
#   - insert a new column E ("constraint"), pushing the old "calculation"
#     column (E) to F
#   - insert a new row 8 for the "age" question, pushing the trailing rows
#     (old 8..15) down to 9..16
#   - fill in the new header/cells
#   - update the sheet's selection to match the post-edit cursor position

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("survey")

# Insert the new "constraint" column between D (relevant) and the old E
# (calculation); this shifts the calculation column + its data to F.
$ws.Columns.Item(5).Insert()

# Insert the new row for the "age" question between "calc" (row 7) and
# "end group" (old row 8); this shifts rows 8..15 down to 9..16.
$ws.Rows.Item(8).Insert()

# New survey row: decimal "age" question with a label and a constraint.
$ws.Range("A8").Value = "decimal"
$ws.Range("B8").Value = "age"
$ws.Range("C8").Value = "Your age:"

# New column header + the constraint formula text for the age row.
$ws.Range("E1").Value = "constraint"
$ws.Range("E8").Value = ". < 150"

# Narrow the new constraint column (the old-E/calculation width carries
# over to the new F column automatically via the insert above).
$ws.Range("E1").EntireColumn.ColumnWidth = 9

# Match the workbook's post-edit cursor position.
$ws.Range("F21").Select()
